$d = $word.ActiveDocument

function Replace-InParagraph {
    param(
        [int]$Index,
        [string]$OldText,
        [string]$NewText
    )
    # Locate the text with Find (no in-place Replace parameter) and then set
    # Range.Text directly - Find's built-in Replace pass runs text through
    # smart-quote AutoFormat (straight ' -> curly '); a plain Range.Text
    # assignment does not, so it round-trips characters like apostrophes
    # byte-for-byte.
    $p = $d.Paragraphs.Item($Index)
    $searchRange = $p.Range.Duplicate
    $ok = $searchRange.Find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false)
    if (-not $ok) {
        Write-Host "WARNING: replace failed at paragraph $Index for text: $OldText"
        return
    }
    $searchRange.Text = $NewText
}

function Insert-ParagraphAfter {
    param(
        [int]$Index,
        [string]$NewText
    )
    $p = $d.Paragraphs.Item($Index)
    $p.Range.InsertParagraphAfter()
    $newPara = $d.Paragraphs.Item($Index + 1)
    $newPara.Range.Text = $NewText
}

# --- PARTNER - Siege Analytics bullets ---
Replace-InParagraph 10 "• Conduct comprehensive quantitative and qualitative research studies using Python, R, SPSS, and Stata for political candidates and organizations" "• Lead comprehensive polling and research studies for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in strategic spending decisions"
Replace-InParagraph 11 "• Architect cloud-based data warehouse solutions on AWS (EC2, RDS, S3) processing billions of records for electoral analytics" "• Architect enterprise-scale cloud data warehouse solutions on AWS (EC2, RDS, S3) processing millions of records with millions of columns for electoral analytics and demographic analysis"
Replace-InParagraph 12 "• Design scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets" "• Design and implement scalable ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial and demographic datasets"
Replace-InParagraph 13 "• Develop custom analytical tools and algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering" "• Develop advanced analytical tools and machine learning algorithms using Python, Pandas, NumPy, and Scikit-learn for fraud detection and spatial clustering"
Replace-InParagraph 14 "• Manage complex client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications" "• Manage strategic client relationships across political, nonprofit, and technology sectors using Django/GeoDjango web applications"
Replace-InParagraph 15 "• Lead technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices" "• Drive technical architecture decisions for data-intensive applications using Docker, Git, and modern DevOps practices"

# --- DATA PRODUCTS MANAGER - Helm/Murmuration bullets ---
Replace-InParagraph 18 "• Conceived and developed framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES" "• Conceived and developed comprehensive data framework using Python, Pandas, and PostgreSQL to clean, validate, and normalize government data from Census, BLS, and NCES"
Replace-InParagraph 19 "• Built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions" "• Architected and built multi-tenant data warehouse and data lake using Snowflake, dbt, and AWS processing millions of records with millions of columns for longitudinal analysis across attitudinal, behavioral, demographic, economic and geographical dimensions"
Replace-InParagraph 20 "• Trained analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI" "• Led training initiatives for analytical and engineering staff on open source geospatial technology (QGIS, GRASS, OSGeo) for analysis, segmentation, and visualization using Tableau and PowerBI"
Replace-InParagraph 21 "• Wrote five-year strategic plans for developing data warehouse using Scala, PySpark, and Apache Spark that became basis of company's distinguishing products" "• Developed five-year strategic plans for data warehouse architecture using Scala, PySpark, and Apache Spark that became foundation of company's distinguishing products"
Replace-InParagraph 22 "• Managed teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices" "• Led cross-functional teams of seven to eleven engineers, designers, analysts, and external stakeholders using Agile methodologies and modern DevOps practices"

# --- SENIOR ANALYST - Myers Research bullet ---
Replace-InParagraph 34 "• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research" "• Provided strategic counsel to Democratic campaigns, political actors, and NGOs through quantitative and qualitative research affecting millions of dollars in campaign spending decisions"

# --- RESEARCH DIRECTOR - Progressive Change Campaign Committee bullets ---
Replace-InParagraph 37 "• Managed critical research operations for political campaigns" "• Conceived, architected, and engineered FLEEM web application using Twilio API for thousands of simultaneous phone calls"
Replace-InParagraph 38 "• Conducted comprehensive polling and demographic analysis" "• Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren"
Replace-InParagraph 39 "• Developed strategic recommendations based on data analysis" "• Built tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"
Replace-InParagraph 40 "• Led research team in support of progressive political initiatives" "• Designed survey deployment system facilitating thousands of simultaneous phone surveys, saving PAC nearly `$1 million annually in polling costs"

# --- PROGRAMMER - Lake Research Partners section heading + bullets ---
Replace-InParagraph 54 "Political Research and Data Analysis" "Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns"
Replace-InParagraph 55 "• Developed data analysis tools for political polling and research" "• Built the first collaborative and multi-actor contributed poll of polls used by the Democratic Party"
Replace-InParagraph 56 "• Built statistical models for voter behavior analysis" "• Developed system that later became the Polling Consortium Database at The Analyst Institute"
Replace-InParagraph 57 "• Created data visualization tools for research presentations" "• Worked on all aspects of questionnaire design, sampling, reporting and analysis for Congressional, Senate and Presidential elections affecting millions of dollars in campaign spending decisions"
Replace-InParagraph 58 "• Supported senior researchers with technical analysis and reporting" "• Conducted statistical modeling and analysis using SPSS, ArcGIS, Quantum GIS, GRASS, Stata, OSCAR, PostgreSQL, PostGIS, and Oracle"

# --- FIELD DIRECTOR - The Feldman Group section heading + bullets ---
Replace-InParagraph 60 "Political Field Operations and Data Management" "Political Polling, Focus Groups and Demographic Analysis for Democratic Campaigns"
Replace-InParagraph 61 "• Managed field operations for political campaigns and research projects" "• Administered all quantitative and qualitative research operations for presidential, gubernatorial, congressional, and senatorial campaigns affecting millions of dollars in spending decisions"
Replace-InParagraph 62 "• Developed data collection and management systems for field work" "• Managed team of 6 research analysts and field staff for comprehensive survey fielding at multi-million dollar research firm"
Replace-InParagraph 63 "• Trained field staff on data collection protocols and quality control" "• Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings"
Replace-InParagraph 64 "• Analyzed field data to inform campaign strategy and research findings" "• Created custom reports and data visualizations based on specific client requirements"

# --- Insertions (process from bottom-most index to top-most so indices stay valid) ---
# After original paragraph 64 ("• Analyzed field data..." -> "• Created custom reports...")
Insert-ParagraphAfter 64 "• Introduced mapping and geospatial analysis into standard reporting procedures"
Insert-ParagraphAfter 65 "• Enhanced value of research deliverables through advanced analytical techniques using SPSS, OSCAR, PHP, and MySQL"

# After original paragraph 58 ("• Supported senior researchers..." -> "• Conducted statistical modeling...")
Insert-ParagraphAfter 58 "• Pioneered integration of advanced mapping techniques into standard reports including choropleths and hexagonal grid maps"
Insert-ParagraphAfter 59 "• Developed innovative approaches to visualizing demographic and market data for enhanced client understanding"

# After original paragraph 40 ("• Led research team..." -> "• Designed survey deployment system...")
Insert-ParagraphAfter 40 "• Managed comprehensive research operations for progressive political initiatives and candidates"

Write-Host "Done. Final paragraph count:" $d.Paragraphs.Count
